# 19th April 1st update
# Fill in column AL (18/04/2020) values for the affected states/rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AL3").Value = 2
$ws.Range("AL7").Value = 1
$ws.Range("AL8").Value = 2
$ws.Range("AL10").Value = 186
$ws.Range("AL12").Value = 277
$ws.Range("AL13").Value = 9
$ws.Range("AL14").Value = 1
$ws.Range("AL15").Value = 13
$ws.Range("AL16").Value = 1
$ws.Range("AL17").Value = 25
$ws.Range("AL18").Value = 4
$ws.Range("AL20").Value = 92
$ws.Range("AL21").Value = 328
$ws.Range("AL26").Value = 1
$ws.Range("AL28").Value = 23
$ws.Range("AL29").Value = 122
$ws.Range("AL30").Value = 49
$ws.Range("AL31").Value = 43
$ws.Range("AL33").Value = 125
$ws.Range("AL34").Value = 2
